# Thêm 5 testcase Nhật Ký trong TestCase.xlsx
# Adds 5 new test-case rows (17-21) to Sheet1, matching the style of the
# existing "Giải Trí" rows (14-16), resizes Table1 / AutoFilter / the sheet
# dimension to A1:H21, and moves the selection/scroll position to the end
# of the newly-added data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Clone the row formatting for the 5 new rows from existing rows that
#    already carry the right style set.
#    Row 16 -> styles 20/21/21/21/21/21/22/23 (used by rows 17,19,20,21)
#    Row 13 -> styles 12/5/5/5/5/5/(none)/14   (used by row 18)
# ---------------------------------------------------------------------
$ws.Range("A16:H16").Copy() | Out-Null
$ws.Range("A17:H17").PasteSpecial(-4122) | Out-Null
$ws.Range("A16:H16").Copy() | Out-Null
$ws.Range("A19:H19").PasteSpecial(-4122) | Out-Null
$ws.Range("A16:H16").Copy() | Out-Null
$ws.Range("A20:H20").PasteSpecial(-4122) | Out-Null
$ws.Range("A16:H16").Copy() | Out-Null
$ws.Range("A21:H21").PasteSpecial(-4122) | Out-Null

$ws.Range("A13:F13").Copy() | Out-Null
$ws.Range("A18:F18").PasteSpecial(-4122) | Out-Null
$ws.Range("H13").Copy() | Out-Null
$ws.Range("H18").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Fill in the cell values, in the exact order the strings were first
#    authored, so the shared-string table grows the same way.
# ---------------------------------------------------------------------

# -- Row 17 : FUNC_NhatKy_TC1 ------------------------------------------------
$ws.Cells.Item(17, 1).Value = "FUNC_NhatKy_TC1"
$ws.Cells.Item(17, 2).Value = "Chọn Giải Trí Nhật Ký"
$ws.Cells.Item(17, 6).Value = "Khi chọn chức năng giải trí từ màn hình chính . Màn hình sẽ hiện ra các trò chơi giải trí cho mình lựa chọn: nhật ký , game , đố vui. Chọn Nhật Ký. Hiện ra màn hình có liên quan đến Nhật Ký"

# -- Row 18 : FUNC_NhatKy_TC2 ------------------------------------------------
$ws.Cells.Item(18, 1).Value = "FUNC_NhatKy_TC2"
$ws.Cells.Item(18, 2).Value = "Viết nhật ký"

$ws.Cells.Item(17, 5).Value = "Chọn  Giải Trí từ màn hình chính. Chọn Nhật Ký. "

$ws.Cells.Item(18, 6).Value = "Khi chọn chức năng giải trí màn hình sẽ hiện ra các trò chơi giải trí để lựa chọn : chơi game , nhật ký , đố vui. Chọn Nhật ký . Màn hình hiện ra cửa số cho mình ghi pass. Nếu ghi pass đúng thì sẽ hiện ra màn hình cho mình ghi nhật ký."

# -- Row 19 : FUNC_NhatKy_TC3 ------------------------------------------------
$ws.Cells.Item(19, 1).Value = "FUNC_NhatKy_TC3"
$ws.Cells.Item(19, 2).Value = "Lưu nhật ký"

$ws.Cells.Item(18, 5).Value = "Chọn chức năng giải trí từ màn hình chính. Chọn Nhật Ký để giải trí . Ghi pass , Viết nhật ký."

$ws.Cells.Item(19, 5).Value = "Chọn chức năng giải trí từ màn hình , chọn nhật ký , ghi pass , viết nhật ký , lưu đoạn nhật ký vừa ghi."

$ws.Cells.Item(19, 6).Value = "Khi chọn chức năng giải trí thì màn hình hiện ra dánh sách các trò giải trí: nhật ký , chơi game , đố vui. Khi chọn nhật ký thì màn hình sẽ hiện ra cửa sổ để ghi pass . Nếu ghi đúng pass , hiện ra màn hình để ghi nhật ký . Trước khi thoát khỏi nhật ký thì có thông báo lưu lại đoạn nhật ký vừa ghi hay không . "

$ws.Cells.Item(20, 6).Value = "Khi chọn chức năng giải trí thì màn hình hiện ra danh sách các trò giải trí: nhật ký , chơi game , đố vui. Khi chọn nhật ký thì màn hình sẽ hiện ra cửa sổ để ghi pass . Nếu ghi đúng pass , không ghi pass , trở lại màn hình trước đó."

$ws.Cells.Item(20, 2).Value = "Trở lại màn hình trước đó"

$ws.Cells.Item(20, 5).Value = "Chọn chức năng giải trí từ màn hình , chọn nhật ký , trở lại màn hình trước đó ."

# -- Row 20 : FUNC_NhatKy_TC4 ------------------------------------------------
$ws.Cells.Item(20, 1).Value = "FUNC_NhatKy_TC4"

# -- Row 21 : FUNC_NhatKy_TC5 ------------------------------------------------
$ws.Cells.Item(21, 1).Value = "FUNC_NhatKy_TC5"

$ws.Cells.Item(21, 5).Value = "Chọn chức năng giải trí từ màn hình , chọn nhật ký , trở lại màn hình chính của phần mềm"

$ws.Cells.Item(21, 6).Value = "Trở lại màn hình chính của phần mềm bất cứ lúc nào."

# -- Columns that reuse already-existing shared strings ---------------------
$ws.Cells.Item(17, 3).Value = "Không có"
$ws.Cells.Item(17, 4).Value = "Không có"
$ws.Cells.Item(18, 3).Value = "Không có"
$ws.Cells.Item(18, 4).Value = "Không có"
$ws.Cells.Item(19, 3).Value = "Không có"
$ws.Cells.Item(19, 4).Value = "Không có"
$ws.Cells.Item(20, 3).Value = "Không có"
$ws.Cells.Item(20, 4).Value = "Không có"
$ws.Cells.Item(21, 3).Value = "Không có"
$ws.Cells.Item(21, 4).Value = "Không có"
$ws.Cells.Item(21, 2).Value = "Trở lại màn hình chính của phần mềm"

# ---------------------------------------------------------------------
# 3. Row heights, matching the source row each new row was cloned from.
# ---------------------------------------------------------------------
$ws.Rows.Item(17).RowHeight = 51
$ws.Rows.Item(18).RowHeight = 63.75
$ws.Rows.Item(19).RowHeight = 76.5
$ws.Rows.Item(20).RowHeight = 63.75
$ws.Rows.Item(21).RowHeight = 38.25

# ---------------------------------------------------------------------
# 4. Grow the table / autofilter / sheet dimension to cover the new rows.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:H21")) | Out-Null

# ---------------------------------------------------------------------
# 5. Move the view so the new rows are visible, mirroring the author's
#    final cursor position.
# ---------------------------------------------------------------------
$ws.Range("B21").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
